$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'285.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'2.75%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'28.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'4.18%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.063"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'4.17%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.06479"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.94%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'7.229"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'4.00%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.344"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'14.06%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9126"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'4.18%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1543"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.01%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.06493"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'26.57%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07663"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'3.43%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.02984"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'1.03%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.08946"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.39%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001593"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'1.95%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0006546"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'2.74%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006103"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.14%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.458"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.61%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.367"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'1.83%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'-1.39%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'1.33%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1342"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-0.46%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.971"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'1.84%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D24").Value = "'0.04462"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.91%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001179"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.10%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004322"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'11.79%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E29").Value = "'-15.67%"
$ws.Range("E29").Style = "Normal"
$ws.Range("D40").Value = "'0.04147"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-0.36%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006778"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-0.63%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1231"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'4.96%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002102"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'4.08%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'3.84%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005399"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'1.48%"
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "BOLO"
$ws.Range("C46").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D46").Value = "'1.933"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'14.73%"
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "CoinbaseStockToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D47").Value = "'0.01851"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.05%"
$ws.Range("E47").Style = "Normal"
